$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-11 Wednesday" "2025-06-12 Thursday"

Replace-Text "588×2=1176" "364×9=3276"
Replace-Text "965×4=3860" "569×3=1707"
Replace-Text "653×6=3918" "764×7=5348"
Replace-Text "414×5=2070" "151×4=604"
Replace-Text "353×8=2824" "954×3=2862"
Replace-Text "472×9=4248" "933×8=7464"
Replace-Text "668×6=4008" "846×4=3384"
Replace-Text "824×8=6592" "400×5=2000"
Replace-Text "558×2=1116" "379×6=2274"
Replace-Text "539×6=3234" "693×8=5544"
Replace-Text "917×6=5502" "843×9=7587"
Replace-Text "917×9=8253" "257×6=1542"
Replace-Text "860×8=6880" "942×7=6594"
Replace-Text "602×8=4816" "569×3=1707"
Replace-Text "880×8=7040" "759×3=2277"
Replace-Text "241×3=723" "642×2=1284"
Replace-Text "745×8=5960" "499×6=2994"
Replace-Text "101×5=505" "578×3=1734"
Replace-Text "720×4=2880" "249×4=996"
Replace-Text "755×6=4530" "937×8=7496"
Replace-Text "691×6=4146" "554×8=4432"
Replace-Text "509×8=4072" "136×2=272"
Replace-Text "465×7=3255" "882×9=7938"
Replace-Text "342×2=684" "529×6=3174"
Replace-Text "860×9=7740" "724×7=5068"
